# Fix a typo in the "hint" note on Sheet1 (I3): "as sheets auxiliares..."
# should read "nas sheets auxiliares..." and leave the active selection on E2,
# matching where the author was working when the file was saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("I3").Value = "nas sheets auxiliares, os valores de busca tem que estar ordenados"

$ws.Activate() | Out-Null
$ws.Range("E2").Select() | Out-Null
